$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the commit diff.
# Cells that look like plain numbers must be forced to Text (matching the
# original inlineStr cell type) via a temporary "@" number format, then the
# cell style is reset back to Normal so no stray style index is left behind.

$ws.Range("D2").Value = "30.525.34"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.687.88"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0624"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "1.931.52"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.68%  "
$ws.Range("D14").Value = "1.701.07"
$ws.Range("E14").Value = "  +4.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.620"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("D17").Value = "30.559.68"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0501"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").Value = "1.508.41"
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.586"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.27%  "
$ws.Range("D48").Value = "1.822.64"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("D51").Value = "0.0₆0116"
$ws.Range("E51").Value = "  +7.47%  "
